$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '64.511.16'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.441.02'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -0.65%  '

# Row 4
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '572.84'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -0.55%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '158.70'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -1.65%  '

# Row 7
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '3.440.07'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -0.75%  '

# Row 9
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -5.52%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.21'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -0.18%  '

# Row 11
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -2.46%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.440'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -1.16%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.034.63'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -0.69%  '

# Row 14
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -0.11%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '27.46'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -2.78%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.0000173'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -9.47%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '64.560.84'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.452.53'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -0.49%  '

# Row 19
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -4.33%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.73'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -4.21%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '379.73'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '7.95'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -2.16%  '

# Row 23
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '72.27'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -0.16%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.528'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -3.88%  '

# Row 26
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.93'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.71%  '

# Row 28
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +0.26%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.989'
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -1.22%  '

# Row 30
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -5.08%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.07'
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.97%  '

# Row 32
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -1.74%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '23.21'
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -1.57%  '

# Row 34
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -1.86%  '

# Row 35
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -2.97%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '161.15'
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -0.15%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.87'
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -3.24%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.880.15'
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -3.07%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0744'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -3.93%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '26.24'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -1.34%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.794'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +2.05%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '4.52'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -0.53%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '42.90'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '6.48'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -3.15%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0310'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -3.09%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '25.76'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '

# Row 47
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +11.90%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '321.15'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +2.45%  '

# Row 49
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -2.89%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '6.45'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -2.76%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.841'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -2.86%  '

